$d = $word.ActiveDocument

function Replace-Text($searchText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find.Execute could not find: $searchText"
    }
}

# 1) "... and Bulk Data File Registration APIs to optionally ..."
#    -> "... and Bulk Data File Registration REST APIs to optionally ..."
Replace-Text "Bulk Data File Registration APIs " "Bulk Data File Registration REST APIs "

# 2) "(Presently, the parent collection can be created, but the metadata of an
#    existing parent cannot be updated)" ->
#    "Presently, the metadata can be specified if the parent collection is
#    being created, but the metadata of an existing parent cannot be updated)"
Replace-Text `
    "(Presently, the parent collection can be created, but the metadata of an existing parent cannot be updated)" `
    "Presently, the metadata can be specified if the parent collection is being created, but the metadata of an existing parent cannot be updated)"

# 3) " Removed the system properties from the " ->
#    " Removed the system specific properties from the "
Replace-Text " Removed the system properties from the " " Removed the system specific properties from the "

# 4) "Updated help information in the DME Web Application and CLU interface" ->
#    "Updated help information in the DME Web Application and in the CLU interface"
Replace-Text `
    "Updated help information in the DME Web Application and CLU interface" `
    "Updated help information in the DME Web Application and in the CLU interface"

Write-Host "All 4 replacements applied successfully."
